$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1100.6666
$ws.Range("I18").Value = 1020.8
$ws.Range("K18").Value = 1020.8
$ws.Range("M18").Value = -736.8
$ws.Range("H53").Value = 391.08334
$ws.Range("I53").Value = 655.5
$ws.Range("J53").Value = 126.666664
$ws.Range("K53").Value = 655.5
$ws.Range("L53").Value = 126.666664
$ws.Range("M53").Value = -18.5
$ws.Range("N53").Value = -1400.666664

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()
$ws.Range("H122").Value = 2210.889
$ws.Range("I122").Value = 2103.3103
$ws.Range("J122").Value = 2656.5715
$ws.Range("K122").Value = 6309.9309
$ws.Range("L122").Value = 7969.7145
$ws.Range("M122").Value = -3859.9309
$ws.Range("N122").Value = -12869.7145
$ws.Range("H132").Value = 2375.375
$ws.Range("I132").Value = 2132.7441
$ws.Range("K132").Value = 6398.2323
$ws.Range("M132").Value = -3868.2323

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2874
$ws.Range("I86").Value = 2525.7778
$ws.Range("J86").Value = 3222.2222
$ws.Range("K86").Value = 2525.7778
$ws.Range("L86").Value = 3222.2222
$ws.Range("M86").Value = -1402.7778
$ws.Range("N86").Value = -5468.2222
$ws.Range("H89").Value = 2874
$ws.Range("I89").Value = 2525.7778
$ws.Range("J89").Value = 3222.2222
$ws.Range("K89").Value = 12628.889
$ws.Range("L89").Value = 16111.111
$ws.Range("M89").Value = -7012.888999999999
$ws.Range("N89").Value = -27343.111
$ws.Range("H134").Value = 2199.762
$ws.Range("I134").Value = 2169.75
$ws.Range("K134").Value = 6509.25
$ws.Range("M134").Value = -3974.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 213.76923
$ws.Range("I7").Value = 249.9
$ws.Range("K7").Value = 249.9
$ws.Range("M7").Value = -136.9
$ws.Range("H31").Value = 23451.883
$ws.Range("I31").Value = 34483.938
$ws.Range("K31").Value = 34483.938
$ws.Range("M31").Value = -34188.938
$ws.Range("H34").Value = 23451.883
$ws.Range("I34").Value = 34483.938
$ws.Range("K34").Value = 34483.938
$ws.Range("M34").Value = -34281.938
$ws.Range("H37").Value = 19993
$ws.Range("I37").Value = 14990
$ws.Range("J37").Value = 29999
$ws.Range("K37").Value = 14990
$ws.Range("L37").Value = 29999
$ws.Range("M37").Value = -14883
$ws.Range("N37").Value = -30213
$ws.Range("H58").Value = 2808.111
$ws.Range("I58").Value = 2669.932
$ws.Range("J58").Value = 8888
$ws.Range("K58").Value = 2669.932
$ws.Range("L58").Value = 8888
$ws.Range("M58").Value = -2466.932
$ws.Range("N58").Value = -9294
$ws.Range("H62").Value = 5400.5
$ws.Range("I62").Value = 4201.8
$ws.Range("J62").Value = 6599.2
$ws.Range("K62").Value = 4201.8
$ws.Range("L62").Value = 6599.2
$ws.Range("M62").Value = -3577.8
$ws.Range("N62").Value = -7847.2
$ws.Range("H65").Value = 5400.5
$ws.Range("I65").Value = 4201.8
$ws.Range("J65").Value = 6599.2
$ws.Range("K65").Value = 21009
$ws.Range("L65").Value = 32996
$ws.Range("M65").Value = -17889
$ws.Range("N65").Value = -39236
$ws.Range("H132").Value = 3919
$ws.Range("I132").Value = 3875
$ws.Range("K132").Value = 11625
$ws.Range("M132").Value = -9095
$ws.Range("H134").Value = 17769.5
$ws.Range("I134").Value = 8531.467000000001
$ws.Range("K134").Value = 25594.401
$ws.Range("M134").Value = -23059.401
$ws.Range("H136").Value = 2808.111
$ws.Range("I136").Value = 2669.932
$ws.Range("J136").Value = 8888
$ws.Range("K136").Value = 8009.795999999999
$ws.Range("L136").Value = 26664
$ws.Range("M136").Value = -5459.795999999999
$ws.Range("N136").Value = -31764

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 536.3333
$ws.Range("I2").Value = 536.3333
$ws.Range("K2").Value = 3217.9998
$ws.Range("M2").Value = -3104.9998
$ws.Range("H39").Value = 1674.5834
$ws.Range("J39").Value = 1674.5834
$ws.Range("L39").Value = 5023.7502
$ws.Range("N39").Value = -5611.7502
$ws.Range("H50").Value = 1003.8571
$ws.Range("J50").Value = 1532.375
$ws.Range("L50").Value = 4597.125
$ws.Range("N50").Value = -5559.125
$ws.Range("H53").Value = 1003.8571
$ws.Range("J53").Value = 1532.375
$ws.Range("L53").Value = 4597.125
$ws.Range("N53").Value = -5559.125
$ws.Range("H62").Value = 4918.4443
$ws.Range("I62").Value = 3253.2
$ws.Range("J62").Value = 7000
$ws.Range("K62").Value = 9759.599999999999
$ws.Range("L62").Value = 21000
$ws.Range("M62").Value = -9073.599999999999
$ws.Range("N62").Value = -22372
$ws.Range("H63").Value = 403.83334
$ws.Range("I63").Value = 431
$ws.Range("J63").Value = 349.5
$ws.Range("K63").Value = 1293
$ws.Range("L63").Value = 1048.5
$ws.Range("M63").Value = -544
$ws.Range("N63").Value = -2546.5
$ws.Range("H65").Value = 4918.4443
$ws.Range("I65").Value = 3253.2
$ws.Range("J65").Value = 7000
$ws.Range("K65").Value = 29278.8
$ws.Range("L65").Value = 63000
$ws.Range("M65").Value = -25846.8
$ws.Range("N65").Value = -69864
$ws.Range("H66").Value = 403.83334
$ws.Range("I66").Value = 431
$ws.Range("J66").Value = 349.5
$ws.Range("K66").Value = 3879
$ws.Range("L66").Value = 3145.5
$ws.Range("M66").Value = -135
$ws.Range("N66").Value = -10633.5
$ws.Range("H80").Value = 3997.5
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 3997.5
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 11992.5
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -13864.5
$ws.Range("H83").Value = 3997.5
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 3997.5
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 35977.5
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -45337.5
$ws.Range("H122").Value = 2465.6667
$ws.Range("J122").Value = 2366.6667
$ws.Range("L122").Value = 21300.0003
$ws.Range("N122").Value = -26200.0003
$ws.Range("H132").Value = 1213.0465
$ws.Range("I132").Value = 1170.0264
$ws.Range("J132").Value = 1540
$ws.Range("K132").Value = 10530.2376
$ws.Range("L132").Value = 13860
$ws.Range("M132").Value = -8000.2376
$ws.Range("N132").Value = -18920

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 34234
$ws.Range("J12").Value = 1352
$ws.Range("L12").Value = 1352
$ws.Range("N12").Value = -1632
$ws.Range("H38").Value = 23000
$ws.Range("J38").Value = 23000
$ws.Range("L38").Value = 23000
$ws.Range("N38").Value = -23926
$ws.Range("H57").Value = 10863.182
$ws.Range("J57").Value = 14999
$ws.Range("L57").Value = 14999
$ws.Range("N57").Value = -16639
$ws.Range("H132").Value = 4444.7827
$ws.Range("I132").Value = 4546.7
$ws.Range("J132").Value = 3765.3333
$ws.Range("K132").Value = 13640.1
$ws.Range("L132").Value = 11295.9999
$ws.Range("M132").Value = -11110.1
$ws.Range("N132").Value = -16355.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("M20").ClearContents()
$ws.Range("N20").ClearContents()
$ws.Range("H22").Value = 2682.0454
$ws.Range("I22").Value = 1603.1666
$ws.Range("J22").Value = 3086.625
$ws.Range("K22").Value = 1603.1666
$ws.Range("L22").Value = 3086.625
$ws.Range("M22").Value = -1308.1666
$ws.Range("N22").Value = -3676.625
$ws.Range("H27").Value = 2682.0454
$ws.Range("I27").Value = 1603.1666
$ws.Range("J27").Value = 3086.625
$ws.Range("K27").Value = 1603.1666
$ws.Range("L27").Value = 3086.625
$ws.Range("M27").Value = -1496.1666
$ws.Range("N27").Value = -3300.625
$ws.Range("H40").Value = 5767.6924
$ws.Range("I40").Value = 5348
$ws.Range("K40").Value = 5348
$ws.Range("M40").Value = -5212
$ws.Range("H100").Value = 7499.4443
$ws.Range("I100").Value = 2199.2
$ws.Range("K100").Value = 2199.2
$ws.Range("M100").Value = -1658.2
$ws.Range("H122").Value = 83340500
$ws.Range("I122").Value = 166672000
$ws.Range("J122").Value = 8999.666999999999
$ws.Range("K122").Value = 500016000
$ws.Range("L122").Value = 26999.001
$ws.Range("M122").Value = -500013550
$ws.Range("N122").Value = -31899.001
$ws.Range("H132").Value = 4475.5713
$ws.Range("I132").Value = 4299.3335
$ws.Range("K132").Value = 12898.0005
$ws.Range("M132").Value = -10368.0005
$ws.Range("H136").Value = 3286.8708
$ws.Range("I136").Value = 2356.4783
$ws.Range("J136").Value = 5961.75
$ws.Range("K136").Value = 7069.4349
$ws.Range("L136").Value = 17885.25
$ws.Range("M136").Value = -4519.4349
$ws.Range("N136").Value = -22985.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 4038.4614
$ws.Range("I29").Value = 4038.4614
$ws.Range("K29").Value = 4038.4614
$ws.Range("M29").Value = -3748.4614
$ws.Range("H122").Value = 4000.8
$ws.Range("I122").Value = 3668
$ws.Range("J122").Value = 4500
$ws.Range("K122").Value = 11004
$ws.Range("L122").Value = 13500
$ws.Range("M122").Value = -8554
$ws.Range("N122").Value = -18400
$ws.Range("H132").Value = 1432.5
$ws.Range("I132").Value = 1432.5
$ws.Range("K132").Value = 4297.5
$ws.Range("M132").Value = -1767.5
$ws.Range("H136").Value = 2167.6428
$ws.Range("I136").Value = 2121
$ws.Range("J136").Value = 2299.0908
$ws.Range("K136").Value = 6363
$ws.Range("L136").Value = 6897.2724
$ws.Range("M136").Value = -3813
$ws.Range("N136").Value = -11997.2724
